$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "D" = 0.127
    "E" = 0.191
    "I" = 0.00007382190070637282
    "J" = 0.00006134495974191544
    "K" = 53
    "L" = 0.4372937293729373
    "M" = 0.011
    "N" = 0.00003006285870456409
    "O" = 0.0002075471698113207
    "P" = 0.011
    "Q" = 0.00003006285870456409
    "R" = 0.0002075471698113207
    "U" = 387.2
    "V" = 1.058212626400656
    "W" = 0.1620299602568022
    "X" = 0.05441067181199521
    "Y" = 0.107619288444807
    "Z" = 0.4129257142472143
    "AA" = 0.00002533091131689704
    "AB" = 0.04237760892340605
    "AC" = -0.04235227801208916
    "AD" = 258.3
    "AE" = 0.01526392817193807
    "AF" = 258.3152639281719
    "AG" = -128.8847360718281
    "AH" = 0.4138240104905479
    "AI" = 0.3886103989874454
    "AJ" = -0.5437824296028753
    "AK" = -0.4644239536502999
    "AN" = 21525
    "AP" = -10740.39467265234
}

foreach ($row in 2, 3) {
    foreach ($col in $changes.Keys) {
        $ws.Range("$col$row").Value = $changes[$col]
    }
}
